$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.166.21'
$ws.Cells.Item(2, 5).Value = '  -2.12%  '

$ws.Cells.Item(3, 4).Value = '2.576.56'
$ws.Cells.Item(3, 5).Value = '  -2.43%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).Value = '''517.74'
$ws.Cells.Item(5, 5).Value = '  -2.44%  '

$ws.Cells.Item(6, 4).Value = '''138.46'
$ws.Cells.Item(6, 5).Value = '  -5.26%  '

$ws.Cells.Item(7, 4).Value = '''0.998'
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 4).Value = '''0.563'
$ws.Cells.Item(8, 5).Value = '  -1.09%  '

$ws.Cells.Item(9, 4).Value = '2.591.52'
$ws.Cells.Item(9, 5).Value = '  -2.63%  '

$ws.Cells.Item(10, 4).Value = '''6.49'
$ws.Cells.Item(10, 5).Value = '  -2.70%  '

$ws.Cells.Item(11, 4).Value = '''0.0995'
$ws.Cells.Item(11, 5).Value = '  -4.38%  '

$ws.Cells.Item(12, 5).Value = '  -2.55%  '

$ws.Cells.Item(13, 5).Value = '  +0.87%  '

$ws.Cells.Item(14, 4).Value = '3.032.83'
$ws.Cells.Item(14, 5).Value = '  -2.60%  '

$ws.Cells.Item(15, 4).Value = '58.206.85'
$ws.Cells.Item(15, 5).Value = '  -2.08%  '

$ws.Cells.Item(16, 4).Value = '''20.26'
$ws.Cells.Item(16, 5).Value = '  -2.64%  '

$ws.Cells.Item(17, 4).Value = '2.583.99'
$ws.Cells.Item(17, 5).Value = '  -3.62%  '

$ws.Cells.Item(18, 4).Value = '''0.0000131'
$ws.Cells.Item(18, 5).Value = '  -3.70%  '

$ws.Cells.Item(19, 4).Value = '''336.29'
$ws.Cells.Item(19, 5).Value = '  -2.14%  '

$ws.Cells.Item(20, 4).Value = '''4.29'
$ws.Cells.Item(20, 5).Value = '  -3.26%  '

$ws.Cells.Item(21, 4).Value = '''10.11'
$ws.Cells.Item(21, 5).Value = '  -4.61%  '

$ws.Cells.Item(22, 4).Value = '''6.36'
$ws.Cells.Item(22, 5).Value = '  -0.34%  '

$ws.Cells.Item(23, 4).Value = '''0.998'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '

$ws.Cells.Item(24, 4).Value = '''66.11'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '

$ws.Cells.Item(25, 5).Value = '  -2.05%  '

$ws.Cells.Item(26, 2).Value = 'Polygon'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(26, 4).Value = '''0.400'
$ws.Cells.Item(26, 5).Value = '  -3.74%  '

$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).Value = '''0.995'
$ws.Cells.Item(27, 5).Value = '  -0.17%  '

$ws.Cells.Item(28, 4).Value = '''6.98'
$ws.Cells.Item(28, 5).Value = '  -2.95%  '

$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$ws.Cells.Item(30, 5).Value = '  -11.09%  '

$ws.Cells.Item(31, 4).Value = '''5.88'
$ws.Cells.Item(31, 5).Value = '  -7.78%  '

$ws.Cells.Item(32, 5).Value = '  -1.69%  '

$ws.Cells.Item(33, 5).Value = '  -3.85%  '

$ws.Cells.Item(34, 5).Value = '  -0.78%  '

$ws.Cells.Item(35, 5).Value = '  -6.65%  '

$ws.Cells.Item(36, 4).Value = '''1.12'
$ws.Cells.Item(36, 5).Value = '  -5.65%  '

$ws.Cells.Item(37, 4).Value = '''36.14'
$ws.Cells.Item(37, 5).Value = '  -0.88%  '

$ws.Cells.Item(38, 2).Value = 'SuiNetwork'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(38, 5).Value = '  -3.36%  '

$ws.Cells.Item(39, 2).Value = 'Fetch.AI'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(39, 4).Value = '''0.832'
$ws.Cells.Item(39, 5).Value = '  -2.29%  '

$ws.Cells.Item(40, 5).Value = '  -2.99%  '

$ws.Cells.Item(41, 4).Value = '''3.49'
$ws.Cells.Item(41, 5).Value = '  -3.66%  '

$ws.Cells.Item(42, 4).Value = '''0.996'
$ws.Cells.Item(42, 5).Value = '  -0.23%  '

$ws.Cells.Item(43, 4).Value = '''271.16'
$ws.Cells.Item(43, 5).Value = '  +0.42%  '

$ws.Cells.Item(44, 4).Value = '''10.71'
$ws.Cells.Item(44, 5).Value = '  -0.09%  '

$ws.Cells.Item(45, 4).Value = '''0.588'
$ws.Cells.Item(45, 5).Value = '  -2.19%  '

$ws.Cells.Item(46, 4).Value = '''0.0943'
$ws.Cells.Item(46, 5).Value = '  -3.69%  '

$ws.Cells.Item(47, 4).Value = '''0.0515'
$ws.Cells.Item(47, 5).Value = '  -3.67%  '

$ws.Cells.Item(48, 4).Value = '''18.37'
$ws.Cells.Item(48, 5).Value = '  -5.21%  '

$ws.Cells.Item(49, 4).Value = '1.970.54'
$ws.Cells.Item(49, 5).Value = '  -3.30%  '

$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50, 4).Value = '''0.0218'
$ws.Cells.Item(50, 5).Value = '  -5.03%  '

$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51, 4).Value = '''4.53'
$ws.Cells.Item(51, 5).Value = '  -4.70%  '
